$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value  = "*maa://22880 (70.2), maa://20276 (82.73), *maa://22749 (62.5)"
$ws.Range("S3").Value  = "maa://24617 (88.35), **maa://20790 (43.94), ***maa://37170 (20.0)"
$ws.Range("W4").Value  = "**maa://32495 (47.93), ***maa://31785 (16.51), ***maa://36683 (26.67)"
$ws.Range("G6").Value  = "maa://24370 (96.15)"
$ws.Range("W7").Value  = "maa://22399 (94.62), *maa://22758 (70.59)"
$ws.Range("C10").Value = "***maa://25695 (19.3), **maa://32237 (37.84), ***maa://34206 (14.29), ***maa://39951 (21.74), **maa://39243 (40.0)"
$ws.Range("S10").Value = "maa://27395 (96.64), maa://22755 (87.5), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("W10").Value = "maa://22301 (97.37), maa://22726 (100.0)"
$ws.Range("W11").Value = "maa://36713 (97.82)"
$ws.Range("W12").Value = "maa://22753 (91.84), *maa://21485 (76.56), maa://37962 (81.25)"
$ws.Range("C13").Value = "maa://24999 (91.38), maa://36673 (91.8), maa://25001 (85.51)"
$ws.Range("W13").Value = "*maa://34957 (76.74), *maa://22768 (53.33)"
$ws.Range("C15").Value = "*maa://22743 (76.19), maa://22734 (83.33), *maa://30808 (63.64), ***maa://36048 (13.33)"
$ws.Range("S16").Value = "maa://22729 (95.17), *maa://28648 (69.09), *maa://36674 (76.92)"
$ws.Range("AE16").Value = "*maa://23911 (61.54), maa://27755 (91.78)"
$ws.Range("K18").Value = "maa://22466 (88.37), *maa://22732 (52.5)"
$ws.Range("AE21").Value = "maa://22524 (94.29), *maa://22432 (74.07)"
$ws.Range("W22").Value = "maa://21282 (98.82), *maa://37649 (64.71)"
$ws.Range("K23").Value = "maa://39756 (92.0), maa://39875 (95.56)"
$ws.Range("W24").Value = "maa://23504 (93.1), maa://29988 (85.93), **maa://22892 (40.14), *maa://25141 (76.86), *maa://36663 (80.0), ***maa://22815 (23.08)"
$ws.Range("W28").Value = "maa://39929 (85.79), ***maa://39723 (15.15), maa://41749 (81.82)"
$ws.Range("K29").Value = "maa://28432 (93.45), *maa://28440 (72.84), maa://31400 (100.0), *maa://28650 (66.67)"
$ws.Range("O29").Value = "*maa://23168 (55.77), **maa://30050 (40.0)"
$ws.Range("K31").Value = "maa://35926 (93.48), maa://36258 (80.52)"
$ws.Range("S32").Value = "maa://41108 (88.0), maa://41238 (94.12)"
$ws.Range("K35").Value = "maa://41296 (97.37)"
$ws.Range("G44").Value = "maa://29768 (97.5), maa://27728 (96.0)"
$ws.Range("G60").Value = "**maa://40438 (40.0)"
